# Update column F ("dSF") values for a set of rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = 0
    14 = 1
    19 = -1
    25 = -4
    28 = 3
    31 = 8
    39 = -11
    40 = -4
    52 = -5
    54 = -5
    58 = -2
    59 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
